# Relabel the "18-21 mos" dataset from "21mos" to "18mos":
#  - rename the worksheet tab itself
#  - on the "accuracy" sheet, retag the corresponding rows (Dataset column)
#    and correct their Date column to the (earlier) 18mos capture date
#  - leave the "accuracy" sheet as the active/selected tab with the
#    matching selection/scroll state, as it was left after the edit

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "21mos" sheet to "18mos" (3rd tab) ---
$sheet18 = $wb.Worksheets.Item(3)
$sheet18.Name = "18mos"

# --- 2. Fix up the "accuracy" sheet rows that belonged to the 21mos dataset ---
$acc = $wb.Worksheets.Item(8)

# Rows 62-93 are the block tagged with the old "21mos" dataset label;
# relabel them "18mos" and normalize their Date (column B) to 43334
# (2018-08-22), replacing the old scattered 21mos-era dates.
$acc.Range("A62:A93").Value = "18mos"
$acc.Range("B62:B93").Value = 43334

# --- 3. Leave "accuracy" as the active sheet/tab with its new selection ---
$acc.Activate()
$acc.Range("J83").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 60
$win.ScrollColumn = 1

# --- 4. Page setup tweak recorded alongside the edit ---
$acc.PageSetup.Orientation = 1
